# Insert 3 new data rows (weekly "Fruta"/"hortaliza" update) right before the
# existing row 410, shifting the old rows 410:480 down to 413:483, matching
# the target dimension A1:T483.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("410:412").Insert()

# Common (unchanged) field values shared by all three new rows, copied from
# the surrounding "Terminal Hortofrutícola Agro Chillán" / Frutilla block.
$mercadoId = 7
$mercado = "Terminal Hortofrutícola Agro Chillán"
$region = "Ñuble"
$codreg = 16
$tipo = "Fruta"
$productoId = 100101
$producto = "Berries"
$categoriaId = 100112025
$categoria = "Frutilla"
$variedad = "Sin especificar"
$unidad = "`$/caja 7 kilos"
$origen = "Provincia de Diguillín"
$kgUnidad = 7

# Row 410: Especial
$ws.Cells.Item(410, 1).Value = $mercadoId
$ws.Cells.Item(410, 2).Value = $mercado
$ws.Cells.Item(410, 3).Value = $region
$ws.Cells.Item(410, 4).Value = 44995
$ws.Cells.Item(410, 5).Value = $codreg
$ws.Cells.Item(410, 6).Value = $tipo
$ws.Cells.Item(410, 7).Value = $productoId
$ws.Cells.Item(410, 8).Value = $producto
$ws.Cells.Item(410, 9).Value = $categoriaId
$ws.Cells.Item(410, 10).Value = $categoria
$ws.Cells.Item(410, 11).Value = $variedad
$ws.Cells.Item(410, 12).Value = "Especial"
$ws.Cells.Item(410, 13).Value = 60
$ws.Cells.Item(410, 14).Value = 7500
$ws.Cells.Item(410, 15).Value = 7500
$ws.Cells.Item(410, 16).Value = 7500
$ws.Cells.Item(410, 17).Value = $unidad
$ws.Cells.Item(410, 18).Value = $origen
$ws.Cells.Item(410, 19).Value = 1071
$ws.Cells.Item(410, 20).Value = $kgUnidad

# Row 411: Primera
$ws.Cells.Item(411, 1).Value = $mercadoId
$ws.Cells.Item(411, 2).Value = $mercado
$ws.Cells.Item(411, 3).Value = $region
$ws.Cells.Item(411, 4).Value = 44995
$ws.Cells.Item(411, 5).Value = $codreg
$ws.Cells.Item(411, 6).Value = $tipo
$ws.Cells.Item(411, 7).Value = $productoId
$ws.Cells.Item(411, 8).Value = $producto
$ws.Cells.Item(411, 9).Value = $categoriaId
$ws.Cells.Item(411, 10).Value = $categoria
$ws.Cells.Item(411, 11).Value = $variedad
$ws.Cells.Item(411, 12).Value = "Primera"
$ws.Cells.Item(411, 13).Value = 60
$ws.Cells.Item(411, 14).Value = 6500
$ws.Cells.Item(411, 15).Value = 6500
$ws.Cells.Item(411, 16).Value = 6500
$ws.Cells.Item(411, 17).Value = $unidad
$ws.Cells.Item(411, 18).Value = $origen
$ws.Cells.Item(411, 19).Value = 929
$ws.Cells.Item(411, 20).Value = $kgUnidad

# Row 412: Segunda
$ws.Cells.Item(412, 1).Value = $mercadoId
$ws.Cells.Item(412, 2).Value = $mercado
$ws.Cells.Item(412, 3).Value = $region
$ws.Cells.Item(412, 4).Value = 44995
$ws.Cells.Item(412, 5).Value = $codreg
$ws.Cells.Item(412, 6).Value = $tipo
$ws.Cells.Item(412, 7).Value = $productoId
$ws.Cells.Item(412, 8).Value = $producto
$ws.Cells.Item(412, 9).Value = $categoriaId
$ws.Cells.Item(412, 10).Value = $categoria
$ws.Cells.Item(412, 11).Value = $variedad
$ws.Cells.Item(412, 12).Value = "Segunda"
$ws.Cells.Item(412, 13).Value = 60
$ws.Cells.Item(412, 14).Value = 5500
$ws.Cells.Item(412, 15).Value = 5500
$ws.Cells.Item(412, 16).Value = 5500
$ws.Cells.Item(412, 17).Value = $unidad
$ws.Cells.Item(412, 18).Value = $origen
$ws.Cells.Item(412, 19).Value = 786
$ws.Cells.Item(412, 20).Value = $kgUnidad
